$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Template Name" header in B1 to "Template"
$ws.Range("B1").Value = "Template"

# The "X"/"Y"/"Z" rows used to carry incrementing counters continuing from the
# Page rows (4,5,6) - reset them to their own 1,2,3 counter
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 3

# Widen column A (drop the old bestFit auto width) and select B2
$ws.Columns("A").ColumnWidth = 15.83
$ws.Range("B2").Select() | Out-Null
